$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new row at position 5 for the new "disabilities Persons" data
#    row. This shifts the old merged "Source" row (was row 5) down to row 6,
#    and the old "Note" row (was row 6) down to row 7. The freshly inserted
#    row 5 inherits formatting from row 4 above it, and row 6/7 keep their
#    original formatting/merge ranges (auto-adjusted by Excel).
# ---------------------------------------------------------------------------
$ws.Rows(5).Insert()

# Remove the old "Note: ... Data is confidential or unavailable." row - it is
# no longer present in the revised sheet.
$ws.Rows(7).Delete()

# ---------------------------------------------------------------------------
# 2. Row 1 - title. Merge across A1:I1, change the title text, and make the
#    alignment centered + wrapped (on top of the already Bold Arial 11 font).
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "The number of persons with disabilities registered in the Unified database of targeted social assistance program in Sagarejo Municipality"
$ws.Range("A1:I1").Merge()
$ws.Range("A1:I1").HorizontalAlignment = -4108
$ws.Range("A1:I1").VerticalAlignment = -4108
$ws.Range("A1:I1").WrapText = $true
$ws.Rows(1).RowHeight = 51

# ---------------------------------------------------------------------------
# 3. Row 2 - unchanged text, but no longer has a custom row height.
# ---------------------------------------------------------------------------
$ws.Rows(2).RowHeight = 14.5

# ---------------------------------------------------------------------------
# 4. Row 3 - blank header cell font becomes Sylfaen (matches new column 1
#    default font), keeps its existing top border.
# ---------------------------------------------------------------------------
$ws.Range("A3").Font.Name = "Sylfaen"
$ws.Range("A3").Font.Size = 11

# ---------------------------------------------------------------------------
# 5. Row 4 - "family with disabilities Persons" data row.
#    A4 keeps its top border but loses the bottom border.
#    B4:I4 become plain numbers, losing the right-alignment override.
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "family with disabilities Persons "
$ws.Range("A4").Borders.Item(9).LineStyle = 0
$ws.Rows(4).RowHeight = 24.75

$row4Vals = @(756,710,666,714,719,734,760,757)
$row4Cols = @("B","C","D","E","F","G","H","I")
for ($i = 0; $i -lt $row4Cols.Length; $i++) {
    $addr = "$($row4Cols[$i])4"
    $ws.Range($addr).Value = $row4Vals[$i]
    $ws.Range($addr).HorizontalAlignment = 1
    $ws.Range($addr).Borders.Item(8).LineStyle = 0
    $ws.Range($addr).Borders.Item(9).LineStyle = 0
}

# ---------------------------------------------------------------------------
# 6. Row 5 - NEW "disabilities Persons" data row (inherited formatting from
#    row 4 on insert). A5 keeps its bottom border but loses the top border.
#    B5:H5 lose the right-alignment override (keep numFmt). I5 keeps its
#    bottom border but loses the top border and the right-alignment
#    override.
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = "disabilities Persons "
$ws.Range("A5").Borders.Item(8).LineStyle = 0
$ws.Rows(5).RowHeight = 21

$row5Vals = @(837,786,740,787,789,804,834,829)
$row5Cols = @("B","C","D","E","F","G","H","I")
for ($i = 0; $i -lt $row5Cols.Length; $i++) {
    $addr = "$($row5Cols[$i])5"
    $ws.Range($addr).Value = $row5Vals[$i]
    $ws.Range($addr).HorizontalAlignment = 1
}
$ws.Range("I5").Borders.Item(8).LineStyle = 0

# ---------------------------------------------------------------------------
# 7. Row 6 - Source row (text unchanged, carried along by the row insert).
#    The merged range (A6:H6) already followed automatically. A6 loses its
#    top border (B6:H6 keep theirs untouched).
# ---------------------------------------------------------------------------
$ws.Range("A6").Borders.Item(8).LineStyle = 0
$ws.Rows(6).RowHeight = 27.75

# ---------------------------------------------------------------------------
# 8. Column width for column A.
# ---------------------------------------------------------------------------
$ws.Columns(1).ColumnWidth = 20

# ---------------------------------------------------------------------------
# 9. Selection shown in the sheet view.
# ---------------------------------------------------------------------------
$ws.Range("A1:I1").Select()
